# edit.ps1 -- "New crime data collected": refresh the weekly 94th Precinct CompStat figures
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Masthead: bump report Volume/Number and the covered week date range ---
$ws.Range("A8").Value2 = "Volume 30   Number  32"
$ws.Range("C9").Value2 = "Report Covering the Week  8/7/2023  Through  8/13/2023"

# --- Crime complaint table (rows 15-29): refreshed weekly figures ---
# Cells whose type flips from a numeric stat to the "no data" placeholder (text "0" / "***.*")
# need their format reset to General via a format-only paste from a known General-formatted
# cell (C22), since re-typing the placeholder alone would keep the old numeric format.
$ws.Range("D15").NumberFormat = "#,##0"
$ws.Range("D15").Value2 = 2
$ws.Range("E15").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("E15").Value2 = -100
$ws.Range("G15").Value2 = 5
$ws.Range("J15").Value2 = 9
$ws.Range("K15").Value2 = -66.666666666666
$ws.Range("C16").Value2 = 4
$ws.Range("D16").Value2 = 5
$ws.Range("E16").Value2 = -20
$ws.Range("F16").Value2 = 10
$ws.Range("G16").Value2 = 12
$ws.Range("H16").Value2 = -16.666666666666
$ws.Range("I16").Value2 = 66
$ws.Range("J16").Value2 = 86
$ws.Range("K16").Value2 = -23.255813953488
$ws.Range("L16").Value2 = 22.222222222222
$ws.Range("M16").Value2 = -18.518518518518
$ws.Range("N16").Value2 = -82.162162162162
$ws.Range("C17").Value2 = 4
$ws.Range("D17").Value2 = 6
$ws.Range("E17").Value2 = -33.333333333333
$ws.Range("I17").Value2 = 69
$ws.Range("J17").Value2 = 96
$ws.Range("K17").Value2 = -28.125
$ws.Range("L17").Value2 = 35.294117647058
$ws.Range("M17").Value2 = 40.816326530612
$ws.Range("N17").Value2 = -56.875
$ws.Range("D18").Value2 = 3
$ws.Range("E18").Value2 = 33.333333333333
$ws.Range("G18").Value2 = 20
$ws.Range("H18").Value2 = -20
$ws.Range("I18").Value2 = 117
$ws.Range("J18").Value2 = 145
$ws.Range("K18").Value2 = -19.310344827586
$ws.Range("L18").Value2 = -7.874015748031
$ws.Range("M18").Value2 = -24.516129032258
$ws.Range("N18").Value2 = -83.380681818181
$ws.Range("C19").Value2 = 8
$ws.Range("D19").Value2 = 13
$ws.Range("E19").Value2 = -38.461538461538
$ws.Range("F19").Value2 = 50
$ws.Range("G19").Value2 = 54
$ws.Range("H19").Value2 = -7.407407407407
$ws.Range("I19").Value2 = 414
$ws.Range("J19").Value2 = 359
$ws.Range("K19").Value2 = 15.320334261838
$ws.Range("L19").Value2 = 70.37037037037
$ws.Range("M19").Value2 = 149.397590361446
$ws.Range("N19").Value2 = 89.041095890411
$ws.Range("C20").Value2 = 4
$ws.Range("D20").Value2 = 3
$ws.Range("E20").Value2 = 33.333333333333
$ws.Range("F20").Value2 = 18
$ws.Range("G20").Value2 = 20
$ws.Range("H20").Value2 = -10
$ws.Range("I20").Value2 = 104
$ws.Range("J20").Value2 = 108
$ws.Range("K20").Value2 = -3.703703703703
$ws.Range("L20").Value2 = 38.666666666666
$ws.Range("M20").Value2 = 13.043478260869
$ws.Range("N20").Value2 = -81.592920353982
$ws.Range("C21").Value2 = 24
$ws.Range("D21").Value2 = 32
$ws.Range("E21").Value2 = -25
$ws.Range("F21").Value2 = 108
$ws.Range("G21").Value2 = 132
$ws.Range("H21").Value2 = -18.181818181818
$ws.Range("I21").Value2 = 773
$ws.Range("J21").Value2 = 804
$ws.Range("K21").Value2 = -3.855721393034
$ws.Range("L21").Value2 = 39.028776978417
$ws.Range("M21").Value2 = 42.095588235294
$ws.Range("N21").Value2 = -61.921182266009
$ws.Range("L22").Value2 = 16.666666666666
$ws.Range("D23").Value2 = "'0"
$ws.Range("E23").Value2 = "'***.*"
$ws.Range("F23").Value2 = 2
$ws.Range("G23").Value2 = 2
$ws.Range("H23").Value2 = 0
$ws.Range("I23").Value2 = 16
$ws.Range("K23").Value2 = 6.666666666666
$ws.Range("L23").Value2 = 14.285714285714
$ws.Range("M23").Value2 = -15.78947368421
$ws.Range("C24").Value2 = 18
$ws.Range("E24").Value2 = -18.181818181818
$ws.Range("F24").Value2 = 68
$ws.Range("G24").Value2 = 84
$ws.Range("H24").Value2 = -19.047619047619
$ws.Range("I24").Value2 = 572
$ws.Range("J24").Value2 = 628
$ws.Range("K24").Value2 = -8.917197452229
$ws.Range("L24").Value2 = 16.024340770791
$ws.Range("M24").Value2 = 61.581920903954
$ws.Range("C25").Value2 = 9
$ws.Range("D25").Value2 = 4
$ws.Range("E25").Value2 = 125
$ws.Range("F25").Value2 = 21
$ws.Range("G25").Value2 = 24
$ws.Range("H25").Value2 = -12.5
$ws.Range("I25").Value2 = 150
$ws.Range("J25").Value2 = 174
$ws.Range("K25").Value2 = -13.793103448275
$ws.Range("L25").Value2 = 17.1875
$ws.Range("M25").Value2 = 7.913669064748
$ws.Range("D26").NumberFormat = "#,##0"
$ws.Range("D26").Value2 = 2
$ws.Range("E26").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("E26").Value2 = -100
$ws.Range("G26").Value2 = 5
$ws.Range("H26").Value2 = -80
$ws.Range("J26").Value2 = 11
$ws.Range("K26").Value2 = -54.545454545454
$ws.Range("C27").NumberFormat = "#,##0"
$ws.Range("C27").Value2 = 2
$ws.Range("D27").Value2 = "'0"
$ws.Range("E27").Value2 = "'***.*"
$ws.Range("F27").Value2 = 5
$ws.Range("H27").Value2 = 25
$ws.Range("I27").Value2 = 30
$ws.Range("K27").Value2 = 66.666666666666
$ws.Range("L27").Value2 = 57.894736842105
$ws.Range("D28").Value2 = "'0"
$ws.Range("E28").Value2 = "'***.*"
$ws.Range("D29").Value2 = "'0"
$ws.Range("E29").Value2 = "'***.*"

# Normalize format of the newly-text placeholder cells to plain General (style 14)
$ws.Range("C22").Copy()
$ws.Range("D23").PasteSpecial(-4122)
$ws.Range("E23").PasteSpecial(-4122)
$ws.Range("D27").PasteSpecial(-4122)
$ws.Range("E27").PasteSpecial(-4122)
$ws.Range("D28").PasteSpecial(-4122)
$ws.Range("E28").PasteSpecial(-4122)
$ws.Range("D29").PasteSpecial(-4122)
$ws.Range("E29").PasteSpecial(-4122)
$excel.CutCopyMode = $false
